$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.121.00"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "3.612.30"
$ws.Range("E3").Value = "  +3.20%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "603.18"
$ws.Range("E5").Value = "  +1.10%  "
$ws.Range("D6").Value = "196.59"
$ws.Range("E6").Value = "  +0.97%  "
$ws.Range("E7").Value = "  +0.67%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -0.59%  "
$ws.Range("D10").Value = "0.649"
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("D11").Value = "53.90"
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").Value = "0.0000306"
$ws.Range("E12").Value = "  +1.88%  "
$ws.Range("D13").Value = "9.56"
$ws.Range("E13").Value = "  +0.50%  "
$ws.Range("D14").Value = "4.191.08"
$ws.Range("E14").Value = "  +3.38%  "
$ws.Range("D15").Value = "13.19"
$ws.Range("E15").Value = "  +4.94%  "
$ws.Range("D16").Value = "592.31"
$ws.Range("E16").Value = "  -2.59%  "
$ws.Range("D17").Value = "19.23"
$ws.Range("E17").Value = "  +1.38%  "
$ws.Range("D18").Value = "70.264.59"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D19").Value = "3.610.88"
$ws.Range("E19").Value = "  +3.02%  "
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("D21").Value = "0.995"
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("D22").Value = "17.70"
$ws.Range("E22").Value = "  -1.41%  "
$ws.Range("E23").Value = "  +0.85%  "
$ws.Range("D24").Value = "102.16"
$ws.Range("E24").Value = "  -1.92%  "
$ws.Range("E25").Value = "  +1.04%  "
$ws.Range("D26").Value = "3.04"
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("E27").Value = "  -0.98%  "
$ws.Range("D28").Value = "9.59"
$ws.Range("E28").Value = "  -0.98%  "
$ws.Range("D29").Value = "34.03"
$ws.Range("E29").Value = "  +1.33%  "
$ws.Range("D30").Value = "4.80"
$ws.Range("E30").Value = "  +3.96%  "
$ws.Range("D31").Value = "7.16"
$ws.Range("E31").Value = "  +1.34%  "
$ws.Range("D32").Value = "12.32"
$ws.Range("E32").Value = "  -2.22%  "
$ws.Range("E33").Value = "  +1.99%  "
$ws.Range("D34").Value = "63.28"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").Value = "0.0₃0895"
$ws.Range("E35").Value = "  +10.42%  "
$ws.Range("D36").Value = "3.925.67"
$ws.Range("E36").Value = "  +5.24%  "
$ws.Range("D37").Value = "3.13"
$ws.Range("E37").Value = "  +3.06%  "
$ws.Range("D38").Value = "528.82"
$ws.Range("E38").Value = "  +5.84%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").Value = "37.32"
$ws.Range("E40").Value = "  +2.00%  "
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("E42").Value = "  -1.44%  "
$ws.Range("E43").Value = "  -1.08%  "
$ws.Range("D44").Value = "0.0456"
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("E45").Value = "  +1.38%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "3.36"
$ws.Range("E46").Value = "  +1.29%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "0.141"
$ws.Range("E47").Value = "  +1.03%  "
$ws.Range("D48").Value = "8.60"
$ws.Range("E48").Value = "  -1.47%  "
$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").Value = "0.000258"
$ws.Range("E49").Value = "  +6.84%  "
$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("D51").Value = "1.33"
$ws.Range("E51").Value = "  +4.18%  "
